# Updates cryptos list with latest price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.115.53"
$ws.Range("E2").Value = "  +3.58%  "
$ws.Range("D3").Value = "'2.451.75"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'568.85"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("D6").Value = "'167.05"
$ws.Range("E6").Value = "  +4.80%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "'0.176"
$ws.Range("E9").Value = "  +12.62%  "
$ws.Range("D10").Value = "'2.449.79"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("E14").Value = "  +8.07%  "
$ws.Range("D15").Value = "'70.007.69"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("D16").Value = "'2.903.85"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'24.22"
$ws.Range("E17").Value = "  +5.53%  "
$ws.Range("D18").Value = "'2.467.50"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "'10.89"
$ws.Range("E19").Value = "  +5.55%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'341.67"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.15"
$ws.Range("E21").Value = "  +4.98%  "
$ws.Range("D22").Value = "'3.90"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("E23").Value = "  +8.84%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'66.48"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "'3.84"
$ws.Range("E26").Value = "  +6.29%  "
$ws.Range("E27").Value = "  +5.92%  "
$ws.Range("D28").Value = "'2.580.61"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "'0.973"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").Value = "'0.0₃0858"
$ws.Range("E30").Value = "  +6.59%  "
$ws.Range("D31").Value = "'7.39"
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("D32").Value = "'461.13"
$ws.Range("E32").Value = "  +10.08%  "
$ws.Range("D33").Value = "'1.25"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").Value = "'160.56"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +8.36%  "
$ws.Range("D38").Value = "'19.11"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'18.22"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").Value = "'0.304"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("D42").Value = "'1.53"
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("D43").Value = "'4.43"
$ws.Range("E43").Value = "  +4.30%  "
$ws.Range("D44").Value = "'38.10"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").Value = "'2.13"
$ws.Range("E46").Value = "  +6.71%  "
$ws.Range("D47").Value = "'134.53"
$ws.Range("E47").Value = "  +4.68%  "
$ws.Range("D48").Value = "'3.40"
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "'0.492"
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").Value = "'0.564"
$ws.Range("E51").Value = "  +2.02%  "
